# Updated cryptos list values (prices / volume / coin reorderings)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds numeric-looking text (e.g. "1.00", "95.811.20").
# Force text format on the whole Price column first so assignments below
# keep their exact literal formatting instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '95.811.20'
$ws.Range('E2').Value = '  +0.54%  '
$ws.Range('D3').Value = '3.557.29'
$ws.Range('E3').Value = '  -1.05%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '238.85'
$ws.Range('E5').Value = '  +0.44%  '
$ws.Range('D6').Value = '652.78'
$ws.Range('E6').Value = '  +0.48%  '
$ws.Range('D7').Value = '1.60'
$ws.Range('E7').Value = '  +9.82%  '
$ws.Range('D8').Value = '0.402'
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('B9').Value = 'USDC'
$ws.Range('C9').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  +0.08%  '
$ws.Range('B10').Value = 'Cardano'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D10').Value = '1.05'
$ws.Range('E10').Value = '  +5.15%  '
$ws.Range('D11').Value = '3.555.67'
$ws.Range('E11').Value = '  -1.05%  '
$ws.Range('D12').Value = '42.98'
$ws.Range('E12').Value = '  +0.82%  '
$ws.Range('E13').Value = '  +1.62%  '
$ws.Range('D14').Value = '6.37'
$ws.Range('E14').Value = '  +1.17%  '
$ws.Range('D15').Value = '4.220.60'
$ws.Range('E15').Value = '  -1.35%  '
$ws.Range('D16').Value = '95.516.90'
$ws.Range('E16').Value = '  +0.34%  '
$ws.Range('E17').Value = '  +1.73%  '
$ws.Range('D18').Value = '3.542.91'
$ws.Range('E18').Value = '  -1.33%  '
$ws.Range('D19').Value = '7.79'
$ws.Range('E19').Value = '  -1.82%  '
$ws.Range('D20').Value = '12.55'
$ws.Range('E20').Value = '  -1.07%  '
$ws.Range('D21').Value = '17.62'
$ws.Range('E21').Value = '  -1.45%  '
$ws.Range('D22').Value = '0.512'
$ws.Range('E22').Value = '  +5.96%  '
$ws.Range('B23').Value = 'SuiNetwork'
$ws.Range('C23').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D23').Value = '3.38'
$ws.Range('E23').Value = '  -6.34%  '
$ws.Range('B24').Value = 'BitcoinCash'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D24').Value = '500.93'
$ws.Range('E24').Value = '  -1.07%  '
$ws.Range('D25').Value = '6.90'
$ws.Range('E25').Value = '  +5.36%  '
$ws.Range('D26').Value = '0.0000197'
$ws.Range('E26').Value = '  +1.22%  '
$ws.Range('D27').Value = '95.21'
$ws.Range('E27').Value = '  -0.72%  '
$ws.Range('D28').Value = '12.69'
$ws.Range('E28').Value = '  +0.66%  '
$ws.Range('D29').Value = '3.749.08'
$ws.Range('E29').Value = '  -0.82%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = '2.99'
$ws.Range('E30').Value = '  -3.91%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = '0.150'
$ws.Range('E31').Value = '  +8.09%  '
$ws.Range('D32').Value = '11.29'
$ws.Range('E32').Value = '  -0.26%  '
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('E34').Value = '  +2.96%  '
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').Value = '31.21'
$ws.Range('E36').Value = '  -1.34%  '
$ws.Range('D37').Value = '8.76'
$ws.Range('E37').Value = '  +7.92%  '
$ws.Range('D38').Value = '611.36'
$ws.Range('E38').Value = '  +5.90%  '
$ws.Range('D39').Value = '0.560'
$ws.Range('E39').Value = '  +0.81%  '
$ws.Range('D40').Value = '1.61'
$ws.Range('E40').Value = '  +9.21%  '
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('E42').Value = '  +0.24%  '
$ws.Range('D43').Value = '0.899'
$ws.Range('E43').Value = '  -1.98%  '
$ws.Range('D44').Value = '1.81'
$ws.Range('E44').Value = '  +5.35%  '
$ws.Range('E45').Value = '  +0.21%  '
$ws.Range('D46').Value = '23.52'
$ws.Range('E46').Value = '  -0.87%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '33.90'
$ws.Range('E47').Value = '  +0.82%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').Value = '2.26'
$ws.Range('E48').Value = '  +1.38%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').Value = '0.0417'
$ws.Range('E49').Value = '  +1.59%  '
$ws.Range('E50').Value = '  +1.12%  '
$ws.Range('E51').Value = '  +0.96%  '
